$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new task "View current lane schedule for a given day" to Adonis's
# task column (E) and to the burn-down list (A).
$ws.Range("E8").Value = "View current lane schedule for a given day"
$ws.Range("A36").Value = "View current lane schedule for a given day"

# Rename the header/persona columns E, F, G from the generic role names to
# the actual persona names.
$ws.Range("E1").Value = "Adonis"
$ws.Range("F1").Value = "Pythagoras"
$ws.Range("G1").Value = "Jørgen"

# Move the selection / view down to the newly added burn-down row.
$ws.Range("A36").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 4
